$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial for rows 2-32; bump each
# from 45602 to 45603 (one day later), matching the diff.
$ws.Range("C2:C32").Value = 45603
